# Auto-generated Excel COM-interop script to apply Sheets diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3131.2778
$ws.Range("I15").Value = 3131.2778
$ws.Range("K15").Value = 9393.8334
$ws.Range("M15").Value = -9224.8334
$ws.Range("H40").Value = 2990.6667
$ws.Range("J40").Value = 2993.5
$ws.Range("L40").Value = 2993.5
$ws.Range("N40").Value = -3343.5
$ws.Range("H53").Value = 2500
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 2500
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 2500
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -3774
$ws.Range("H86").Value = 3290.7273
$ws.Range("I86").Value = 3098.75
$ws.Range("K86").Value = 3098.75
$ws.Range("M86").Value = -1975.75
$ws.Range("H89").Value = 3290.7273
$ws.Range("I89").Value = 3098.75
$ws.Range("K89").Value = 15493.75
$ws.Range("M89").Value = -9877.75
$ws.Range("H115").Value = 461.66666
$ws.Range("J115").Value = 500
$ws.Range("L115").Value = 1500
$ws.Range("N115").Value = -4634
$ws.Range("H132").Value = 6766.731
$ws.Range("I132").Value = 7592.381
$ws.Range("J132").Value = 3299
$ws.Range("K132").Value = 22777.143
$ws.Range("L132").Value = 9897
$ws.Range("M132").Value = -20247.143
$ws.Range("N132").Value = -14957
$ws.Range("H138").Value = 2702.9714
$ws.Range("I138").Value = 2427.2307
$ws.Range("K138").Value = 7281.6921
$ws.Range("M138").Value = -2141.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 241666
$ws.Range("J64").Value = 249999
$ws.Range("L64").Value = 249999
$ws.Range("N64").Value = -250495
$ws.Range("H67").Value = 241666
$ws.Range("J67").Value = 249999
$ws.Range("L67").Value = 249999
$ws.Range("N67").Value = -251715
$ws.Range("H97").Value = 1603.0555
$ws.Range("I97").Value = 985.6
$ws.Range("J97").Value = 2374.875
$ws.Range("K97").Value = 985.6
$ws.Range("L97").Value = 2374.875
$ws.Range("M97").Value = -489.6
$ws.Range("N97").Value = -3366.875
$ws.Range("H122").Value = 2326.65
$ws.Range("I122").Value = 2019.3846
$ws.Range("J122").Value = 2897.2856
$ws.Range("K122").Value = 6058.1538
$ws.Range("L122").Value = 8691.856800000001
$ws.Range("M122").Value = -3608.1538
$ws.Range("N122").Value = -13591.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4590.2856
$ws.Range("J99").Value = 6411.75
$ws.Range("L99").Value = 6411.75
$ws.Range("N99").Value = -9407.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 65361.688
$ws.Range("J58").Value = 4197.8
$ws.Range("L58").Value = 4197.8
$ws.Range("N58").Value = -4603.8
$ws.Range("H86").Value = 7249
$ws.Range("I86").Value = 7249
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 7249
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -6126
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 7249
$ws.Range("I89").Value = 7249
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 36245
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -30629
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 2726.1785
$ws.Range("I107").Value = 1077.2727
$ws.Range("J107").Value = 3793.1177
$ws.Range("K107").Value = 1077.2727
$ws.Range("L107").Value = 3793.1177
$ws.Range("M107").Value = 842.7273
$ws.Range("N107").Value = -7633.1177
$ws.Range("H132").Value = 1933.8
$ws.Range("I132").Value = 1941.25
$ws.Range("J132").Value = 1904
$ws.Range("K132").Value = 5823.75
$ws.Range("L132").Value = 5712
$ws.Range("M132").Value = -3293.75
$ws.Range("N132").Value = -10772
$ws.Range("H136").Value = 65361.688
$ws.Range("J136").Value = 4197.8
$ws.Range("L136").Value = 12593.4
$ws.Range("N136").Value = -17693.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 91.76922999999999
$ws.Range("I12").Value = 84
$ws.Range("K12").Value = 252
$ws.Range("M12").Value = -79

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 12500
$ws.Range("J34").Value = 12500
$ws.Range("L34").Value = 12500
$ws.Range("N34").Value = -13036
$ws.Range("H41").Value = 24618.5
$ws.Range("J41").Value = 29990
$ws.Range("L41").Value = 29990
$ws.Range("N41").Value = -30700
$ws.Range("H52").Value = 10030
$ws.Range("I52").Value = 10030
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 10030
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -9771
$ws.Range("N52").ClearContents()
$ws.Range("H76").Value = 12500
$ws.Range("J76").Value = 12500
$ws.Range("L76").Value = 12500
$ws.Range("N76").Value = -13130
$ws.Range("H79").Value = 12500
$ws.Range("J79").Value = 12500
$ws.Range("L79").Value = 12500
$ws.Range("N79").Value = -14684
$ws.Range("H80").Value = 10150
$ws.Range("I80").Value = 3533.3333
$ws.Range("K80").Value = 3533.3333
$ws.Range("M80").Value = -2535.3333
$ws.Range("H83").Value = 10150
$ws.Range("I83").Value = 3533.3333
$ws.Range("K83").Value = 17666.6665
$ws.Range("M83").Value = -12674.6665
$ws.Range("H97").Value = 1362.2
$ws.Range("I97").Value = 1330.7273
$ws.Range("K97").Value = 1330.7273
$ws.Range("M97").Value = -834.7273
$ws.Range("H107").Value = 63061.188
$ws.Range("I107").Value = 91034.55
$ws.Range("K107").Value = 91034.55
$ws.Range("M107").Value = -89114.55
$ws.Range("H122").Value = 3205.6843
$ws.Range("I122").Value = 2541.6667
$ws.Range("K122").Value = 7625.000100000001
$ws.Range("M122").Value = -5175.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4186.6
$ws.Range("I7").Value = 2733.25
$ws.Range("K7").Value = 2733.25
$ws.Range("M7").Value = -2621.25
$ws.Range("H16").Value = 2016.6364
$ws.Range("I16").Value = 2131.5557
$ws.Range("K16").Value = 2131.5557
$ws.Range("M16").Value = -1961.5557
$ws.Range("H45").Value = 44608
$ws.Range("J45").Value = 48333.332
$ws.Range("L45").Value = 48333.332
$ws.Range("N45").Value = -49147.332
$ws.Range("H55").Value = 2044.2858
$ws.Range("I55").Value = 545
$ws.Range("J55").Value = 2644
$ws.Range("K55").Value = 545
$ws.Range("L55").Value = 2644
$ws.Range("M55").Value = -372
$ws.Range("N55").Value = -2990
$ws.Range("H82").Value = 2964.2307
$ws.Range("I82").Value = 3007.8333
$ws.Range("J82").Value = 2926.8572
$ws.Range("K82").Value = 3007.8333
$ws.Range("L82").Value = 2926.8572
$ws.Range("M82").Value = -2646.8333
$ws.Range("N82").Value = -3648.8572
$ws.Range("H85").Value = 2964.2307
$ws.Range("I85").Value = 3007.8333
$ws.Range("J85").Value = 2926.8572
$ws.Range("K85").Value = 3007.8333
$ws.Range("L85").Value = 2926.8572
$ws.Range("M85").Value = -1759.8333
$ws.Range("N85").Value = -5422.8572
$ws.Range("H126").Value = 4186.6
$ws.Range("I126").Value = 2733.25
$ws.Range("K126").Value = 8199.75
$ws.Range("M126").Value = -5729.75
$ws.Range("H130").Value = 37567
$ws.Range("J130").Value = 37567
$ws.Range("L130").Value = 37567
$ws.Range("N130").Value = -47607

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19684
$ws.Range("J41").Value = 19684
$ws.Range("L41").Value = 19684
$ws.Range("N41").Value = -20464
$ws.Range("H58").Value = 17749.25
$ws.Range("J58").Value = 15997
$ws.Range("L58").Value = 15997
$ws.Range("N58").Value = -16613
$ws.Range("H100").Value = 975.94116
$ws.Range("I100").Value = 828.0714
$ws.Range("J100").Value = 1666
$ws.Range("K100").Value = 1656.1428
$ws.Range("L100").Value = 3332
$ws.Range("M100").Value = -1115.1428
$ws.Range("N100").Value = -4414
$ws.Range("H107").Value = 2049.25
$ws.Range("I107").Value = 1119.4
$ws.Range("K107").Value = 3358.2
$ws.Range("M107").Value = -1438.2
$ws.Range("H122").Value = 3024.75
$ws.Range("I122").Value = 2650
$ws.Range("K122").Value = 7950
$ws.Range("M122").Value = -5500
